# Plastics_Database.xlsx edit: normalize the "Size" values on the
# Template sheet (dropping "corti"/micro-sign/space/comma variants down
# to the canonical spellings already used elsewhere in the workbook),
# add the two brand-new XL size strings to the Training Lists lookup
# column, and wire up a dropdown for the new lookup row.

$wb  = $excel.ActiveWorkbook
$tpl = $wb.Worksheets.Item("Template")
$lst = $wb.Worksheets.Item("Training Lists")

# --- Template!B2:B24 -- normalize Size text -------------------------------
$tpl.Range("B2").Value  = "1000ul"
$tpl.Range("B3").Value  = "20ul"
$tpl.Range("B4").Value  = "10ul"
$tpl.Range("B5").Value  = "1000ul"
$tpl.Range("B6").Value  = "1000 ul XL"
$tpl.Range("B7").Value  = "200ul"
$tpl.Range("B8").Value  = "20ul"
$tpl.Range("B9").Value  = "10/20 ul XL"
$tpl.Range("B10").Value = "10ul"
$tpl.Range("B11").Value = "20ul"
$tpl.Range("B12").Value = "200ul"
$tpl.Range("B13").Value = "1000ul"
$tpl.Range("B14").Value = "50mL"
$tpl.Range("B15").Value = "15mL"
$tpl.Range("B16").Value = "0,2mL"
$tpl.Range("B17").Value = "0,2mL"
$tpl.Range("B18").Value = "1,5mL"
$tpl.Range("B19").Value = "1,5mL"
$tpl.Range("B20").Value = "1,5ml"
$tpl.Range("B21").Value = "2mL"
$tpl.Range("B22").Value = "1,5mL"
$tpl.Range("B23").Value = "2ml"
$tpl.Range("B24").Value = "5mL"

# --- Training Lists!B18:B19 -- add the new lookup entries ----------------
$lst.Range("B18").Value = "1000 ul XL"
$lst.Range("B19").Value = "10/20 ul XL"
$lst.Range("B19").HorizontalAlignment = -4108
$lst.Range("B19").VerticalAlignment = -4108

# New validation dropdown on B19, sourced from the growing B column list.
$lst.Range("B19").Validation.Add(3, 1, 1, "=`$B`$2:`$B`$1048576")

# --- selections -------------------------------------------------------
[void]$lst.Activate()
[void]$lst.Range("B20").Select()

[void]$tpl.Activate()
[void]$tpl.Range("H37").Select()
